$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Duplicate formatting of existing rows 2 and 3 onto new rows 12 and 13
$ws.Range("A2:F2").Copy()
$ws.Range("A12:F12").PasteSpecial(-4122)
$ws.Range("A3:F3").Copy()
$ws.Range("A13:F13").PasteSpecial(-4122)

# Row 12: 2026-01-06, 四方坪站
$ws.Range("A12").Value = 46028
$ws.Range("B12").Value = "四方坪站"
$ws.Range("C12").Value = 13867.3
$ws.Range("D12").Value = 9790.75
$ws.Range("E12").Value = 3084.97
$ws.Range("F12").Value = 634

# Row 13: 2026-01-06, 高岭站
$ws.Range("A13").Value = 46028
$ws.Range("B13").Value = "高岭站"
$ws.Range("C13").Value = 6538
$ws.Range("D13").Value = 5436.53
$ws.Range("E13").Value = 1804.24
$ws.Range("F13").Value = 217

$ws.Range("I14").Select()
